# Add auto generate xlsx function
# Append two new rows of translated greetings to Sheet1, right below the
# existing "이것은 / 테스트 / 이니라" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "안녕"

$ws.Range("A4").Value = "하이"
$ws.Range("B4").Value = "헬로"
$ws.Range("C4").Value = "봉쥬르"

# Mirror the author's final selection position (C5) recorded in the file.
$ws.Range("C5").Select()
